$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark every test case's "Run" column (D) as "Y", running all Search test
# cases. Rows 45, 118 and 119 are already "Y" and stay untouched.
for ($r = 2; $r -le 127; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value() -ne "Y") {
        $cell.Value = "Y"
    }
}

# Scroll back to the top and select the full results column, matching the
# final view state after running all the test cases.
$ws.Activate()
$ws.Range("D2:D127").Select()
